# Set "Diferencia Stock" (column L) to 0 for the specified rows,
# and update the "Total_Ajuste_Stock:" sum cell (C63) to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,5,6,7,8,9,13,18,20,21,22,23,25,26,28,30,31,32,34,38,40,41,42,43,45,46)

foreach ($r in $rows) {
    $ws.Range("L$r").Value = 0
}

$ws.Range("C63").Value = 0
